$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario tracker")

$xlPasteFormats = -4122

# Template cells already carrying the three fills we need to reuse
# (copying formats avoids minting brand-new, duplicate style/fill entries)
$greenTemplate = $ws.Range("B3")   # light-green "in db" fill
$redTemplate   = $ws.Range("D23")  # solid-red fill
$noneTemplate  = $ws.Range("D11")  # no-fill cell

# Rows whose status flips to "in db" (light-green fill)
$greenRows = @(7,13,16,27,28)
foreach ($r in $greenRows) {
    $cell = $ws.Range("B$r")
    $greenTemplate.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
    $cell.Value = "in db"
}

# Rows whose status flips to "error" (solid-red fill)
$redRows = @(23,24,25,26)
foreach ($r in $redRows) {
    $cell = $ws.Range("B$r")
    $redTemplate.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
    $cell.Value = "error"
}

# Stray notes next to the carbon-tax / subsidy rows are cleared back to
# plain, unfilled cells
$clearCells = @("C23","D23","C24","C25","C26")
foreach ($addr in $clearCells) {
    $cell = $ws.Range($addr)
    $noneTemplate.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
    $cell.ClearContents()
}

# The stray "will likely error" note next to CarbonStabilityLow is removed
# entirely (cell no longer present in the sheet)
$ws.Range("C27").Clear() | Out-Null

$excel.CutCopyMode = 0

# Selection moves to B26 on the scenario tracker tab
$ws.Range("B26").Select() | Out-Null

# The workbook view no longer pins a custom first visible sheet tab
$wb.Windows.Item(1).DisplayedFirstSheet = 0
